# first_release_qoq_DOMUSE.xlsx edit
# - Rename the "value" header (column B) to "first_release_value"
# - Replace the date/value series with the refreshed vintage (83 rows, 2007-11 .. 2025-08)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: B1 "value" -> "first_release_value" (A1 "date" is unchanged)
$ws.Range("B1").Value = "first_release_value"

# The data block grows from 52 to 83 rows (A1:B53 -> A1:B84). Stretch the existing
# date-column formatting (YYYY-MM-DD style, bold, border -- style of A2) down over
# the whole new range before the values are written so every new row picks it up.
$ws.Range("A2").Copy($ws.Range("A2:A84"))

# Full replacement data for A2:B84 (date serial, value)
$data = New-Object "object[,]" 83,2
$data[0,0] = 38398
$data[0,1] = -1.019024872684525
$data[1,0] = 38487
$data[1,1] = 0.8876024303436765
$data[2,0] = 38579
$data[2,1] = 0.501281665261132
$data[3,0] = 38671
$data[3,1] = 0.783791364788172
$data[4,0] = 38763
$data[4,1] = 0.2322508451440228
$data[5,0] = 38852
$data[5,1] = 0.8665864397470244
$data[6,0] = 38944
$data[6,1] = 0.6793221546917749
$data[7,0] = 39036
$data[7,1] = -0.9228048723025069
$data[8,0] = 39128
$data[8,1] = 1.742611766026243
$data[9,0] = 39217
$data[9,1] = -0.8957516004554691
$data[10,0] = 39309
$data[10,1] = 0.893923595129948
$data[11,0] = 39401
$data[11,1] = -0.2567725410682868
$data[12,0] = 39493
$data[12,1] = 1.910893826230975
$data[13,0] = 39583
$data[13,1] = -0.9623800849065276
$data[14,0] = 39675
$data[14,1] = 1.320073940737448
$data[15,0] = 39767
$data[15,1] = -0.05810994121875979
$data[16,0] = 39859
$data[16,1] = -1.734840982136873
$data[17,0] = 39948
$data[17,1] = -1.4178810011895
$data[18,0] = 40040
$data[18,1] = 1.258487819346726
$data[19,0] = 40132
$data[19,1] = -1.470147873721189
$data[20,0] = 40224
$data[20,1] = 1.38192517089017
$data[21,0] = 40313
$data[21,1] = 1.938408417700344
$data[22,0] = 40405
$data[22,1] = 0.4
$data[23,0] = 40497
$data[23,1] = -0.163634158232469
$data[24,0] = 40589
$data[24,1] = 1.079830393426633
$data[25,0] = 40678
$data[25,1] = 0.2
$data[26,0] = 40770
$data[26,1] = 0.4
$data[27,0] = 40862
$data[27,1] = 0.2303676816657827
$data[28,0] = 40954
$data[28,1] = -0.3
$data[29,0] = 41044
$data[29,1] = -0.4045626487644824
$data[30,0] = 41136
$data[30,1] = 0
$data[31,0] = 41228
$data[31,1] = -0.01847206600469065
$data[32,0] = 41320
$data[32,1] = 0.009235986179263023
$data[33,0] = 41409
$data[33,1] = 0.3863474960573257
$data[34,0] = 41501
$data[34,1] = 0.6964197943645729
$data[35,0] = 41593
$data[35,1] = -0.3183250015702015
$data[36,0] = 41685
$data[36,1] = 1.852186157158073
$data[37,0] = 41774
$data[37,1] = -0.1726927221574073
$data[38,0] = 41866
$data[38,1] = -0.1729977607768376
$data[39,0] = 41958
$data[39,1] = 1.050807574684342
$data[40,0] = 42050
$data[40,1] = 0.5056288600178789
$data[41,0] = 42139
$data[41,1] = -0.2184803162966205
$data[42,0] = 42231
$data[42,1] = 0.7235341094351355
$data[43,0] = 42323
$data[43,1] = 0.8513884674671885
$data[44,0] = 42415
$data[44,1] = 0.8066734233961483
$data[45,0] = 42505
$data[45,1] = -0.1388227614901609
$data[46,0] = 42597
$data[46,1] = 0.4819278240608753
$data[47,0] = 42689
$data[47,1] = 0.6629265129002277
$data[48,0] = 42781
$data[48,1] = 0.1646390629436354
$data[49,0] = 42870
$data[49,1] = 1.164700738417963
$data[50,0] = 42962
$data[50,1] = 0.4317240674915439
$data[51,0] = 43054
$data[51,1] = 0.1359170431485039
$data[52,0] = 43146
$data[52,1] = 0.3857269132374052
$data[53,0] = 43235
$data[53,1] = 0.876336956515118
$data[54,0] = 43327
$data[54,1] = 0.8
$data[55,0] = 43419
$data[55,1] = 0.008724159582257585
$data[56,0] = 43511
$data[56,1] = 0.2
$data[57,0] = 43600
$data[57,1] = 0.4878538807911497
$data[58,0] = 43692
$data[58,1] = -0.4043302599539206
$data[59,0] = 43784
$data[59,1] = 0.7183553771707381
$data[60,0] = 43876
$data[60,1] = -1.5
$data[61,0] = 43966
$data[61,1] = -7.231044133207007
$data[62,0] = 44058
$data[62,1] = 4.729401638091318
$data[63,0] = 44150
$data[63,1] = -0.8905127363963885
$data[64,0] = 44242
$data[64,1] = -2.04269378128221
$data[65,0] = 44331
$data[65,1] = 2.093024636165651
$data[66,0] = 44423
$data[66,1] = 1.098535546956398
$data[67,0] = 44515
$data[67,1] = -0.4717175472572421
$data[68,0] = 44607
$data[68,1] = 0.9401304606753627
$data[69,0] = 44696
$data[69,1] = 0.6392725048137464
$data[70,0] = 44788
$data[70,1] = 0.06357296580725347
$data[71,0] = 44880
$data[71,1] = -1.012166871044968
$data[72,0] = 44972
$data[72,1] = -0.5703626997413522
$data[73,0] = 45061
$data[73,1] = 0.2394101325822788
$data[74,0] = 45153
$data[74,1] = -0.04072131480353391
$data[75,0] = 45245
$data[75,1] = -0.07958838003274593
$data[76,0] = 45337
$data[76,1] = 0.02912383308249389
$data[77,0] = 45427
$data[77,1] = -0.1311265493919933
$data[78,0] = 45519
$data[78,1] = -0.03907468377752821
$data[79,0] = 45611
$data[79,1] = 0.1771324545010202
$data[80,0] = 45703
$data[80,1] = 0.4946531409412387
$data[81,0] = 45792
$data[81,1] = 0.202428137729683
$data[82,0] = 45884
$data[82,1] = 0.208573386070384

$ws.Range("A2:B84").Value = $data

